$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 29, shifting existing rows 29-34 down to 31-36
$ws.Range("A29:A30").EntireRow.Insert()

# New row 29 data (Packham's Triumph, Segunda)
$ws.Range("A29").Value = 1
$ws.Range("B29").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C29").Value = "Arica y Parinacota"
$ws.Range("D29").Value = 44769
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100104
$ws.Range("H29").Value = "Frutos de pepita"
$ws.Range("I29").Value = 100104005
$ws.Range("J29").Value = "Pera"
$ws.Range("K29").Value = "Packham's Triumph"
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 300
$ws.Range("N29").Value = 16000
$ws.Range("O29").Value = 17000
$ws.Range("P29").Value = 16500
$ws.Range("Q29").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 917
$ws.Range("T29").Value = 18

# New row 30 data (Winter Nelis, Segunda)
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value = 44769
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100104
$ws.Range("H30").Value = "Frutos de pepita"
$ws.Range("I30").Value = 100104005
$ws.Range("J30").Value = "Pera"
$ws.Range("K30").Value = "Winter Nelis"
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 15000
$ws.Range("O30").Value = 16000
$ws.Range("P30").Value = 15500
$ws.Range("Q30").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 861
$ws.Range("T30").Value = 18

$ws.Range("D29:D30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
